{"js": "// The document body contains a single 20-row x 5-column table whose cells\n// each hold one arithmetic expression (e.g. \"93-66=27\"). The commit\n// replaces the expression text in every cell with a new expression,\n// keeping the table shape and all run/paragraph formatting untouched.\n//\n// Cell values are addressed by position (row-major, matching\n// Table.values / Table.getRange layout) rather than by searching for the\n// old text, because a handful of the new expressions are identical to an\n// *original* expression found elsewhere in the table (e.g. \"45+18=63\" and\n// \"57+6=63\" each appear once as an old value and once as a new value).\n// A sequential find-and-replace on text would risk matching the wrong\n// cell once an earlier replacement created a duplicate; assigning the\n// whole grid at once avoids that ambiguity entirely.\nconst newValues = [\n  [\"83-38=45\", \"3+69=72\", \"58+19=77\", \"18+75=93\", \"42-35=7\"],\n  [\"69+28=97\", \"7+37=44\", \"84-37=47\", \"83-9=74\", \"81-6=75\"],\n  [\"72-4=68\", \"18+27=45\", \"44+29=73\", \"68+9=77\", \"9+55=64\"],\n  [\"28+47=75\", \"35-26=9\", \"90-45=45\", \"80-8=72\", \"18+59=77\"],\n  [\"8+86=94\", \"23-14=9\", \"38+5=43\", \"42-19=23\", \"93-27=66\"],\n  [\"22-16=6\", \"93-28=65\", \"81-26=55\", \"57+18=75\", \"45+18=63\"],\n  [\"29+46=75\", \"9+77=86\", \"62-17=45\", \"34-18=16\", \"44-28=16\"],\n  [\"92-7=85\", \"90-88=2\", \"76+8=84\", \"18+13=31\", \"77-39=38\"],\n  [\"22+39=61\", \"17+75=92\", \"10-9=1\", \"73-39=34\", \"55-26=29\"],\n  [\"34-27=7\", \"58+18=76\", \"80-35=45\", \"3+29=32\", \"4+48=52\"],\n  [\"53+29=82\", \"53-28=25\", \"7+17=24\", \"83-17=66\", \"54-27=27\"],\n  [\"48+29=77\", \"91-43=48\", \"34+49=83\", \"30-11=19\", \"49+29=78\"],\n  [\"82-78=4\", \"67-48=19\", \"18+39=57\", \"46+28=74\", \"58+5=63\"],\n  [\"50-41=9\", \"70-51=19\", \"65+18=83\", \"85-77=8\", \"4+19=23\"],\n  [\"28+4=32\", \"84-56=28\", \"87-79=8\", \"6+36=42\", \"84-7=77\"],\n  [\"81-17=64\", \"8+9=17\", \"27+47=74\", \"93-67=26\", \"67-19=48\"],\n  [\"9+53=62\", \"93-34=59\", \"29+39=68\", \"35-28=7\", \"90-57=33\"],\n  [\"29+6=35\", \"62-24=38\", \"83-55=28\", \"67+25=92\", \"41-3=38\"],\n  [\"91-48=43\", \"23+28=51\", \"15+27=42\", \"12+49=61\", \"57+6=63\"],\n  [\"39+47=86\", \"22+59=81\", \"63-48=15\", \"77-59=18\", \"93-27=66\"],\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length || table.values[0].length !== newValues[0].length) {\n  throw new Error(\n    \"Table shape mismatch: expected \" + newValues.length + \"x\" + newValues[0].length +\n    \" but found \" + table.rowCount + \"x\" + table.values[0].length\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document body contains a single 20-row x 5-column table whose cells\n# each hold one arithmetic expression (e.g. \"93-66=27\"). The commit replaces\n# the expression text in every cell with a new expression, keeping the table\n# shape and all run/paragraph formatting untouched.\n#\n# Cells are addressed by (row, column) position rather than located via\n# Find/Execute on the old text, because a few of the new expressions equal\n# an *original* expression sitting elsewhere in the table (e.g. \"45+18=63\"\n# and \"57+6=63\" each occur once as an old value and once as a new value).\n# A sequential text search-and-replace could latch onto an already-updated\n# cell instead of the intended one; per-cell positional assignment avoids\n# that ambiguity entirely.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"83-38=45\", \"3+69=72\", \"58+19=77\", \"18+75=93\", \"42-35=7\"),\n    @(\"69+28=97\", \"7+37=44\", \"84-37=47\", \"83-9=74\", \"81-6=75\"),\n    @(\"72-4=68\", \"18+27=45\", \"44+29=73\", \"68+9=77\", \"9+55=64\"),\n    @(\"28+47=75\", \"35-26=9\", \"90-45=45\", \"80-8=72\", \"18+59=77\"),\n    @(\"8+86=94\", \"23-14=9\", \"38+5=43\", \"42-19=23\", \"93-27=66\"),\n    @(\"22-16=6\", \"93-28=65\", \"81-26=55\", \"57+18=75\", \"45+18=63\"),\n    @(\"29+46=75\", \"9+77=86\", \"62-17=45\", \"34-18=16\", \"44-28=16\"),\n    @(\"92-7=85\", \"90-88=2\", \"76+8=84\", \"18+13=31\", \"77-39=38\"),\n    @(\"22+39=61\", \"17+75=92\", \"10-9=1\", \"73-39=34\", \"55-26=29\"),\n    @(\"34-27=7\", \"58+18=76\", \"80-35=45\", \"3+29=32\", \"4+48=52\"),\n    @(\"53+29=82\", \"53-28=25\", \"7+17=24\", \"83-17=66\", \"54-27=27\"),\n    @(\"48+29=77\", \"91-43=48\", \"34+49=83\", \"30-11=19\", \"49+29=78\"),\n    @(\"82-78=4\", \"67-48=19\", \"18+39=57\", \"46+28=74\", \"58+5=63\"),\n    @(\"50-41=9\", \"70-51=19\", \"65+18=83\", \"85-77=8\", \"4+19=23\"),\n    @(\"28+4=32\", \"84-56=28\", \"87-79=8\", \"6+36=42\", \"84-7=77\"),\n    @(\"81-17=64\", \"8+9=17\", \"27+47=74\", \"93-67=26\", \"67-19=48\"),\n    @(\"9+53=62\", \"93-34=59\", \"29+39=68\", \"35-28=7\", \"90-57=33\"),\n    @(\"29+6=35\", \"62-24=38\", \"83-55=28\", \"67+25=92\", \"41-3=38\"),\n    @(\"91-48=43\", \"23+28=51\", \"15+27=42\", \"12+49=61\", \"57+6=63\"),\n    @(\"39+47=86\", \"22+59=81\", \"63-48=15\", \"77-59=18\", \"93-27=66\"),\n)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n        $cell = $table.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $newValues[$r][$c]\n    }\n}\n"}
